$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-03-02 Sunday" "2025-03-03 Monday"

Replace-Text "50×59=" "43×29="
Replace-Text "63×74=" "94×67="
Replace-Text "32×82=" "60×57="
Replace-Text "73×46=" "15×29="
Replace-Text "80×67=" "95×65="
Replace-Text "39×38=" "62×68="
Replace-Text "54×32=" "61×73="
Replace-Text "32×84=" "50×52="
Replace-Text "24×87=" "70×78="
Replace-Text "84×19=" "31×58="
Replace-Text "25×37=" "69×52="
Replace-Text "82×32=" "18×25="
Replace-Text "36×70=" "77×76="
Replace-Text "67×11=" "27×40="
Replace-Text "40×70=" "29×14="
Replace-Text "98×67=" "20×81="
Replace-Text "49×76=" "41×89="
Replace-Text "70×73=" "34×23="
Replace-Text "30×84=" "45×41="
Replace-Text "66×75=" "34×30="
Replace-Text "34×65=" "26×32="
Replace-Text "89×58=" "88×76="
Replace-Text "83×53=" "11×46="
Replace-Text "62×38=" "73×59="
Replace-Text "26×16=" "83×18="
